$d = $word.ActiveDocument
$enDash = [char]8211

# --- 1. Heading: "DATA" / "301  -" / " Project timeline " were split around a
#        stale grammar-check marker ("301  -" flagged gramStart/gramEnd).
#        Re-affirming the full text merges the runs and drops the stale
#        proofing marks (no wording change).
$d.Content.Find.Execute(
    "DATA301  - Project timeline ", $true, $false, $false, $false, $false,
    $true, 1, $false, "DATA301  - Project timeline ", 2) | Out-Null

# --- 2. The actual content edit for this commit: the "Data integrated
#        together" milestone moves from the 15th to the 19th of August.
$d.Content.Find.Execute(
    "Data integrated together: Wednesday 15", $true, $false, $false, $false,
    $false, $true, 1, $false, "Data integrated together: Wednesday 19", 2) | Out-Null

# --- 3. "Has there been an increase ... comparison to 2018." - a stale
#        gramStart/gramEnd pair wraps "2018." which, being the very last
#        piece of text in its paragraph, needs a one-character anchor past
#        the end of the sentence so the whole flagged span (including the
#        trailing marker) gets rebuilt as a single clean run. A throw-away
#        placeholder character is appended, the sentence (plus placeholder)
#        is re-affirmed via Find/Replace, then the placeholder is removed.
$p11 = $d.Paragraphs.Item(11)
$d.Range($p11.Range.End - 1, $p11.Range.End - 1).InsertAfter("#") | Out-Null
$d.Content.Find.Execute(
    "Has there been an increase in Domestic violence reports as a result of COVID-19 lockdown in comparison to 2018.#",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Has there been an increase in Domestic violence reports as a result of COVID-19 lockdown in comparison to 2018.#",
    2) | Out-Null
$p11 = $d.Paragraphs.Item(11)
$d.Range($p11.Range.End - 2, $p11.Range.End - 1).Delete()

# --- 4. "Time and place of crime (all NZ) - Kieryn " - stale spell-check
#        marker around "Kieryn" - re-affirm to merge runs and drop the mark.
$d.Content.Find.Execute(
    "Time and place of crime (all NZ) $enDash Kieryn ", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "Time and place of crime (all NZ) $enDash Kieryn ", 2) | Out-Null

# --- 5. "Start bibliography - Kieryn " - same stale spell-check marker
#        cleanup.
$d.Content.Find.Execute(
    "Start bibliography $enDash Kieryn ", $true, $false, $false, $false,
    $false, $true, 1, $false, "Start bibliography $enDash Kieryn ", 2) | Out-Null

# --- 6. "Create scientific trello based on question discussion mindmap -
#        Abbey " - stale spell-check markers around "trello" and "mindmap".
$d.Content.Find.Execute(
    "Create scientific trello based on question discussion mindmap $enDash Abbey ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Create scientific trello based on question discussion mindmap $enDash Abbey ",
    2) | Out-Null

# Remove the stale "_GoBack" bookmark left over from the last edit session
# (sits right after the "... - Abbey " paragraph above).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
